$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.671.82'
$ws.Range('E2').Value = '  -3.24%  '
$ws.Range('D3').Value = '1.731.96'
$ws.Range('E3').Value = '  -5.72%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '237.90'
$ws.Range('E5').Value = '  -7.82%  '
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4906'
$ws.Range('E7').Value = '  -6.59%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '41.63'
$ws.Range('E8').Value = '  -7.15%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.2423'
$ws.Range('E9').Value = '  -23.00%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.06006'
$ws.Range('E10').Value = '  -11.44%  '
$ws.Range('D11').Value = '1.740.85'
$ws.Range('E11').Value = '  -4.88%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.06707'
$ws.Range('E12').Value = '  -13.39%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '14.73'
$ws.Range('E13').Value = '  -20.92%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '77.17'
$ws.Range('E14').Value = '  -11.93%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.5884'
$ws.Range('E15').Value = '  -24.26%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '4.370'
$ws.Range('E16').Value = '  -12.48%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.000'
$ws.Range('E17').Value = '  -0.01%  '
$ws.Range('E18').Value = '  +0.21%  '
$ws.Range('D19').Value = '25.726.51'
$ws.Range('E19').Value = '  -3.08%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.30'
$ws.Range('E20').Value = '  -18.02%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.000006333'
$ws.Range('E21').Value = '  -19.87%  '
$ws.Range('D22').Value = '1.967.04'
$ws.Range('E22').Value = '  -4.87%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '3.899'
$ws.Range('E23').Value = '  -14.95%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.129'
$ws.Range('E24').Value = '  -13.65%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '7.876'
$ws.Range('E25').Value = '  -15.40%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '135.97'
$ws.Range('E26').Value = '  -4.96%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.849'
$ws.Range('E27').Value = '  -15.92%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.442'
$ws.Range('E28').Value = '  -14.21%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '14.26'
$ws.Range('E29').Value = '  -15.38%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '100.15'
$ws.Range('E30').Value = '  -9.35%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08147'
$ws.Range('E31').Value = '  -6.46%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.632'
$ws.Range('E32').Value = '  -12.71%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.255'
$ws.Range('E33').Value = '  -19.64%  '
$ws.Range('B34').Value = 'Frax'
$ws.Range('C34').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.001'
$ws.Range('E34').Value = '  +0.23%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.04317'
$ws.Range('E35').Value = '  -10.98%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.655'
$ws.Range('E36').Value = '  -7.12%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.018'
$ws.Range('E37').Value = '  -10.31%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.6073'
$ws.Range('E38').Value = '  -16.15%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.766'
$ws.Range('E39').Value = '  -10.39%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.070'
$ws.Range('E40').Value = '  -7.58%  '
$ws.Range('E41').Value = '  +0.11%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.01489'
$ws.Range('E42').Value = '  -13.24%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '101.25'
$ws.Range('E43').Value = '  -7.79%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.7899'
$ws.Range('E44').Value = '  -11.33%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.3812'
$ws.Range('E45').Value = '  -20.50%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '5.140'
$ws.Range('E46').Value = '  -13.01%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '6.068'
$ws.Range('E47').Value = '  -20.43%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.05071'
$ws.Range('E48').Value = '  -12.63%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.239'
$ws.Range('E49').Value = '  -11.89%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '29.49'
$ws.Range('E50').Value = '  -14.80%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '51.76'
$ws.Range('E51').Value = '  -13.08%  '
